# string function case update
# Adds 30 new DQL test-case rows (rows 1758-1787) to Sheet1, covering:
#  - a constant-query case (dqlc1_1741/select 1 from $schema52 limit 1)
#  - mid() boundary tests
#  - subString() boundary / tail-substring tests
#  - pow() with a fractional exponent
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1758
$ws.Cells.Item(1758, 1).Value = 'dqlc1_1757'
$ws.Cells.Item(1758, 2).Value = 'y'
$ws.Cells.Item(1758, 3).Value = '从数据表中查询常量'
$ws.Cells.Item(1758, 4).Value = 'SQLFunctions'
$ws.Cells.Item(1758, 5).NumberFormat = "@"
$ws.Cells.Item(1758, 6).Value = 'schema52'
$ws.Cells.Item(1758, 7).NumberFormat = "@"
$ws.Cells.Item(1758, 8).Value = 'dqlc1_1741'
$ws.Cells.Item(1758, 9).Value = 'select 1 from $schema52 limit 1'
$ws.Cells.Item(1758, 10).Value = 'src/test/resources/io.dingodb.test/testdata/mysqlcases/dql/casegroup1/expectedresult/SQLFuncs/queryc1_1757.csv'
$ws.Cells.Item(1758, 11).Value = 'csv_equals'

# Row 1759
$ws.Cells.Item(1759, 1).Value = 'dqlc1_1758'
$ws.Cells.Item(1759, 2).Value = 'y'
$ws.Cells.Item(1759, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1759, 4).Value = 'String function'
$ws.Cells.Item(1759, 5).Value = 'mid'
$ws.Cells.Item(1759, 6).NumberFormat = "@"
$ws.Cells.Item(1759, 7).NumberFormat = "@"
$ws.Cells.Item(1759, 8).NumberFormat = "@"
$ws.Cells.Item(1759, 9).Value = 'select mid(''abc'',0,3)'
$ws.Cells.Item(1759, 10).NumberFormat = "@"
$ws.Cells.Item(1759, 11).Value = 'string_equals'

# Row 1760
$ws.Cells.Item(1760, 1).Value = 'dqlc1_1759'
$ws.Cells.Item(1760, 2).Value = 'y'
$ws.Cells.Item(1760, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1760, 4).Value = 'String function'
$ws.Cells.Item(1760, 5).Value = 'mid'
$ws.Cells.Item(1760, 6).NumberFormat = "@"
$ws.Cells.Item(1760, 7).NumberFormat = "@"
$ws.Cells.Item(1760, 8).NumberFormat = "@"
$ws.Cells.Item(1760, 9).Value = 'select mid('' abc '',0,3)'
$ws.Cells.Item(1760, 10).NumberFormat = "@"
$ws.Cells.Item(1760, 11).Value = 'string_equals'

# Row 1761
$ws.Cells.Item(1761, 1).Value = 'dqlc1_1760'
$ws.Cells.Item(1761, 2).Value = 'y'
$ws.Cells.Item(1761, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1761, 4).Value = 'String function'
$ws.Cells.Item(1761, 5).Value = 'mid'
$ws.Cells.Item(1761, 6).NumberFormat = "@"
$ws.Cells.Item(1761, 7).NumberFormat = "@"
$ws.Cells.Item(1761, 8).NumberFormat = "@"
$ws.Cells.Item(1761, 9).Value = 'select mid(''0123'',0,1)'
$ws.Cells.Item(1761, 10).NumberFormat = "@"
$ws.Cells.Item(1761, 11).Value = 'string_equals'

# Row 1762
$ws.Cells.Item(1762, 1).Value = 'dqlc1_1761'
$ws.Cells.Item(1762, 2).Value = 'y'
$ws.Cells.Item(1762, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1762, 4).Value = 'String function'
$ws.Cells.Item(1762, 5).Value = 'mid'
$ws.Cells.Item(1762, 6).NumberFormat = "@"
$ws.Cells.Item(1762, 7).NumberFormat = "@"
$ws.Cells.Item(1762, 8).NumberFormat = "@"
$ws.Cells.Item(1762, 9).Value = 'select mid(''abcde'',10,2)'
$ws.Cells.Item(1762, 10).NumberFormat = "@"
$ws.Cells.Item(1762, 11).Value = 'string_equals'

# Row 1763
$ws.Cells.Item(1763, 1).Value = 'dqlc1_1762'
$ws.Cells.Item(1763, 2).Value = 'y'
$ws.Cells.Item(1763, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1763, 4).Value = 'String function'
$ws.Cells.Item(1763, 5).Value = 'mid'
$ws.Cells.Item(1763, 6).NumberFormat = "@"
$ws.Cells.Item(1763, 7).NumberFormat = "@"
$ws.Cells.Item(1763, 8).NumberFormat = "@"
$ws.Cells.Item(1763, 9).Value = 'select mid(''2'',2,1)'
$ws.Cells.Item(1763, 10).NumberFormat = "@"
$ws.Cells.Item(1763, 11).Value = 'string_equals'

# Row 1764
$ws.Cells.Item(1764, 1).Value = 'dqlc1_1763'
$ws.Cells.Item(1764, 2).Value = 'y'
$ws.Cells.Item(1764, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1764, 4).Value = 'String function'
$ws.Cells.Item(1764, 5).Value = 'mid'
$ws.Cells.Item(1764, 6).NumberFormat = "@"
$ws.Cells.Item(1764, 7).NumberFormat = "@"
$ws.Cells.Item(1764, 8).NumberFormat = "@"
$ws.Cells.Item(1764, 9).Value = 'select mid(''abc'',4,1)'
$ws.Cells.Item(1764, 10).NumberFormat = "@"
$ws.Cells.Item(1764, 11).Value = 'string_equals'

# Row 1765
$ws.Cells.Item(1765, 1).Value = 'dqlc1_1764'
$ws.Cells.Item(1765, 2).Value = 'y'
$ws.Cells.Item(1765, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1765, 4).Value = 'String function'
$ws.Cells.Item(1765, 5).Value = 'mid'
$ws.Cells.Item(1765, 6).NumberFormat = "@"
$ws.Cells.Item(1765, 7).NumberFormat = "@"
$ws.Cells.Item(1765, 8).NumberFormat = "@"
$ws.Cells.Item(1765, 9).Value = 'select mid(''abcde'',-10,8)'
$ws.Cells.Item(1765, 10).NumberFormat = "@"
$ws.Cells.Item(1765, 11).Value = 'string_equals'

# Row 1766
$ws.Cells.Item(1766, 1).Value = 'dqlc1_1765'
$ws.Cells.Item(1766, 2).Value = 'y'
$ws.Cells.Item(1766, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1766, 4).Value = 'String function'
$ws.Cells.Item(1766, 5).Value = 'mid'
$ws.Cells.Item(1766, 6).NumberFormat = "@"
$ws.Cells.Item(1766, 7).NumberFormat = "@"
$ws.Cells.Item(1766, 8).NumberFormat = "@"
$ws.Cells.Item(1766, 9).Value = 'select mid(''abcdefg'',7.5,2)'
$ws.Cells.Item(1766, 10).NumberFormat = "@"
$ws.Cells.Item(1766, 11).Value = 'string_equals'

# Row 1767
$ws.Cells.Item(1767, 1).Value = 'dqlc1_1766'
$ws.Cells.Item(1767, 2).Value = 'y'
$ws.Cells.Item(1767, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1767, 4).Value = 'String function'
$ws.Cells.Item(1767, 5).Value = 'mid'
$ws.Cells.Item(1767, 6).NumberFormat = "@"
$ws.Cells.Item(1767, 7).NumberFormat = "@"
$ws.Cells.Item(1767, 8).NumberFormat = "@"
$ws.Cells.Item(1767, 9).Value = 'select mid(123,4,1)'
$ws.Cells.Item(1767, 10).NumberFormat = "@"
$ws.Cells.Item(1767, 11).Value = 'string_equals'

# Row 1768
$ws.Cells.Item(1768, 1).Value = 'dqlc1_1767'
$ws.Cells.Item(1768, 2).Value = 'y'
$ws.Cells.Item(1768, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1768, 4).Value = 'String function'
$ws.Cells.Item(1768, 5).Value = 'mid'
$ws.Cells.Item(1768, 6).NumberFormat = "@"
$ws.Cells.Item(1768, 7).NumberFormat = "@"
$ws.Cells.Item(1768, 8).NumberFormat = "@"
$ws.Cells.Item(1768, 9).Value = 'select mid(123.0,0,1)'
$ws.Cells.Item(1768, 10).NumberFormat = "@"
$ws.Cells.Item(1768, 11).Value = 'string_equals'

# Row 1769
$ws.Cells.Item(1769, 1).Value = 'dqlc1_1768'
$ws.Cells.Item(1769, 2).Value = 'y'
$ws.Cells.Item(1769, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1769, 4).Value = 'String function'
$ws.Cells.Item(1769, 5).Value = 'mid'
$ws.Cells.Item(1769, 6).NumberFormat = "@"
$ws.Cells.Item(1769, 7).NumberFormat = "@"
$ws.Cells.Item(1769, 8).NumberFormat = "@"
$ws.Cells.Item(1769, 9).Value = 'select mid(''abcde'',6)'
$ws.Cells.Item(1769, 10).NumberFormat = "@"
$ws.Cells.Item(1769, 11).Value = 'string_equals'

# Row 1770
$ws.Cells.Item(1770, 1).Value = 'dqlc1_1769'
$ws.Cells.Item(1770, 2).Value = 'y'
$ws.Cells.Item(1770, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1770, 4).Value = 'String function'
$ws.Cells.Item(1770, 5).Value = 'mid'
$ws.Cells.Item(1770, 6).NumberFormat = "@"
$ws.Cells.Item(1770, 7).NumberFormat = "@"
$ws.Cells.Item(1770, 8).NumberFormat = "@"
$ws.Cells.Item(1770, 9).Value = 'select mid(''abcde'',-6)'
$ws.Cells.Item(1770, 10).NumberFormat = "@"
$ws.Cells.Item(1770, 11).Value = 'string_equals'

# Row 1771
$ws.Cells.Item(1771, 1).Value = 'dqlc1_1770'
$ws.Cells.Item(1771, 2).Value = 'y'
$ws.Cells.Item(1771, 3).Value = 'mid函数截取起始值越界'
$ws.Cells.Item(1771, 4).Value = 'String function'
$ws.Cells.Item(1771, 5).Value = 'mid'
$ws.Cells.Item(1771, 6).NumberFormat = "@"
$ws.Cells.Item(1771, 7).NumberFormat = "@"
$ws.Cells.Item(1771, 8).NumberFormat = "@"
$ws.Cells.Item(1771, 9).Value = 'select mid(''abcde'',0)'
$ws.Cells.Item(1771, 10).NumberFormat = "@"
$ws.Cells.Item(1771, 11).Value = 'string_equals'

# Row 1772
$ws.Cells.Item(1772, 1).Value = 'dqlc1_1771'
$ws.Cells.Item(1772, 2).Value = 'y'
$ws.Cells.Item(1772, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1772, 4).Value = 'String function'
$ws.Cells.Item(1772, 5).Value = 'subString'
$ws.Cells.Item(1772, 6).NumberFormat = "@"
$ws.Cells.Item(1772, 7).NumberFormat = "@"
$ws.Cells.Item(1772, 8).NumberFormat = "@"
$ws.Cells.Item(1772, 9).Value = 'select subString(''abcefg'',10,2)'
$ws.Cells.Item(1772, 10).NumberFormat = "@"
$ws.Cells.Item(1772, 11).Value = 'string_equals'

# Row 1773
$ws.Cells.Item(1773, 1).Value = 'dqlc1_1772'
$ws.Cells.Item(1773, 2).Value = 'y'
$ws.Cells.Item(1773, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1773, 4).Value = 'String function'
$ws.Cells.Item(1773, 5).Value = 'subString'
$ws.Cells.Item(1773, 6).NumberFormat = "@"
$ws.Cells.Item(1773, 7).NumberFormat = "@"
$ws.Cells.Item(1773, 8).NumberFormat = "@"
$ws.Cells.Item(1773, 9).Value = 'select subString(''http://www.baidu.com'',22,1)'
$ws.Cells.Item(1773, 10).NumberFormat = "@"
$ws.Cells.Item(1773, 11).Value = 'string_equals'

# Row 1774
$ws.Cells.Item(1774, 1).Value = 'dqlc1_1773'
$ws.Cells.Item(1774, 2).Value = 'y'
$ws.Cells.Item(1774, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1774, 4).Value = 'String function'
$ws.Cells.Item(1774, 5).Value = 'subString'
$ws.Cells.Item(1774, 6).NumberFormat = "@"
$ws.Cells.Item(1774, 7).NumberFormat = "@"
$ws.Cells.Item(1774, 8).NumberFormat = "@"
$ws.Cells.Item(1774, 9).Value = 'select subString(''1234567'',9,3)'
$ws.Cells.Item(1774, 10).NumberFormat = "@"
$ws.Cells.Item(1774, 11).Value = 'string_equals'

# Row 1775
$ws.Cells.Item(1775, 1).Value = 'dqlc1_1774'
$ws.Cells.Item(1775, 2).Value = 'y'
$ws.Cells.Item(1775, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1775, 4).Value = 'String function'
$ws.Cells.Item(1775, 5).Value = 'subString'
$ws.Cells.Item(1775, 6).NumberFormat = "@"
$ws.Cells.Item(1775, 7).NumberFormat = "@"
$ws.Cells.Item(1775, 8).NumberFormat = "@"
$ws.Cells.Item(1775, 9).Value = 'select subString(''abc'',4,1)'
$ws.Cells.Item(1775, 10).NumberFormat = "@"
$ws.Cells.Item(1775, 11).Value = 'string_equals'

# Row 1776
$ws.Cells.Item(1776, 1).Value = 'dqlc1_1775'
$ws.Cells.Item(1776, 2).Value = 'y'
$ws.Cells.Item(1776, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1776, 4).Value = 'String function'
$ws.Cells.Item(1776, 5).Value = 'subString'
$ws.Cells.Item(1776, 6).NumberFormat = "@"
$ws.Cells.Item(1776, 7).NumberFormat = "@"
$ws.Cells.Item(1776, 8).NumberFormat = "@"
$ws.Cells.Item(1776, 9).Value = 'select subString(''abc'',0,3)'
$ws.Cells.Item(1776, 10).NumberFormat = "@"
$ws.Cells.Item(1776, 11).Value = 'string_equals'

# Row 1777
$ws.Cells.Item(1777, 1).Value = 'dqlc1_1776'
$ws.Cells.Item(1777, 2).Value = 'y'
$ws.Cells.Item(1777, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1777, 4).Value = 'String function'
$ws.Cells.Item(1777, 5).Value = 'subString'
$ws.Cells.Item(1777, 6).NumberFormat = "@"
$ws.Cells.Item(1777, 7).NumberFormat = "@"
$ws.Cells.Item(1777, 8).NumberFormat = "@"
$ws.Cells.Item(1777, 9).Value = 'select subString('' abc '',0,3)'
$ws.Cells.Item(1777, 10).NumberFormat = "@"
$ws.Cells.Item(1777, 11).Value = 'string_equals'

# Row 1778
$ws.Cells.Item(1778, 1).Value = 'dqlc1_1777'
$ws.Cells.Item(1778, 2).Value = 'y'
$ws.Cells.Item(1778, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1778, 4).Value = 'String function'
$ws.Cells.Item(1778, 5).Value = 'subString'
$ws.Cells.Item(1778, 6).NumberFormat = "@"
$ws.Cells.Item(1778, 7).NumberFormat = "@"
$ws.Cells.Item(1778, 8).NumberFormat = "@"
$ws.Cells.Item(1778, 9).Value = 'select subString(''0123'',0,1)'
$ws.Cells.Item(1778, 10).NumberFormat = "@"
$ws.Cells.Item(1778, 11).Value = 'string_equals'

# Row 1779
$ws.Cells.Item(1779, 1).Value = 'dqlc1_1778'
$ws.Cells.Item(1779, 2).Value = 'y'
$ws.Cells.Item(1779, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1779, 4).Value = 'String function'
$ws.Cells.Item(1779, 5).Value = 'subString'
$ws.Cells.Item(1779, 6).NumberFormat = "@"
$ws.Cells.Item(1779, 7).NumberFormat = "@"
$ws.Cells.Item(1779, 8).NumberFormat = "@"
$ws.Cells.Item(1779, 9).Value = 'select subString(''abcdef'',-10,1)'
$ws.Cells.Item(1779, 10).NumberFormat = "@"
$ws.Cells.Item(1779, 11).Value = 'string_equals'

# Row 1780
$ws.Cells.Item(1780, 1).Value = 'dqlc1_1779'
$ws.Cells.Item(1780, 2).Value = 'y'
$ws.Cells.Item(1780, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1780, 4).Value = 'String function'
$ws.Cells.Item(1780, 5).Value = 'subString'
$ws.Cells.Item(1780, 6).NumberFormat = "@"
$ws.Cells.Item(1780, 7).NumberFormat = "@"
$ws.Cells.Item(1780, 8).NumberFormat = "@"
$ws.Cells.Item(1780, 9).Value = 'select subString(123,4,1)'
$ws.Cells.Item(1780, 10).NumberFormat = "@"
$ws.Cells.Item(1780, 11).Value = 'string_equals'

# Row 1781
$ws.Cells.Item(1781, 1).Value = 'dqlc1_1780'
$ws.Cells.Item(1781, 2).Value = 'y'
$ws.Cells.Item(1781, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1781, 4).Value = 'String function'
$ws.Cells.Item(1781, 5).Value = 'subString'
$ws.Cells.Item(1781, 6).NumberFormat = "@"
$ws.Cells.Item(1781, 7).NumberFormat = "@"
$ws.Cells.Item(1781, 8).NumberFormat = "@"
$ws.Cells.Item(1781, 9).Value = 'select subString(123.0,0,3)'
$ws.Cells.Item(1781, 10).NumberFormat = "@"
$ws.Cells.Item(1781, 11).Value = 'string_equals'

# Row 1782
$ws.Cells.Item(1782, 1).Value = 'dqlc1_1781'
$ws.Cells.Item(1782, 2).Value = 'y'
$ws.Cells.Item(1782, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1782, 4).Value = 'String function'
$ws.Cells.Item(1782, 5).Value = 'subString'
$ws.Cells.Item(1782, 6).NumberFormat = "@"
$ws.Cells.Item(1782, 7).NumberFormat = "@"
$ws.Cells.Item(1782, 8).NumberFormat = "@"
$ws.Cells.Item(1782, 9).Value = 'select subString(''www.baidu.com'' from 0 for 4)'
$ws.Cells.Item(1782, 10).NumberFormat = "@"
$ws.Cells.Item(1782, 11).Value = 'string_equals'

# Row 1783
$ws.Cells.Item(1783, 1).Value = 'dqlc1_1782'
$ws.Cells.Item(1783, 2).Value = 'y'
$ws.Cells.Item(1783, 3).Value = 'substring函数截取起始值越界'
$ws.Cells.Item(1783, 4).Value = 'String function'
$ws.Cells.Item(1783, 5).Value = 'subString'
$ws.Cells.Item(1783, 6).NumberFormat = "@"
$ws.Cells.Item(1783, 7).NumberFormat = "@"
$ws.Cells.Item(1783, 8).NumberFormat = "@"
$ws.Cells.Item(1783, 9).Value = 'select subString(''www.baidu.com'' from 15 for 4)'
$ws.Cells.Item(1783, 10).NumberFormat = "@"
$ws.Cells.Item(1783, 11).Value = 'string_equals'

# Row 1784
$ws.Cells.Item(1784, 1).Value = 'dqlc1_1783'
$ws.Cells.Item(1784, 2).Value = 'y'
$ws.Cells.Item(1784, 3).Value = 'substring函数从指定起始索引处截取到最后'
$ws.Cells.Item(1784, 4).Value = 'String function'
$ws.Cells.Item(1784, 5).Value = 'subString'
$ws.Cells.Item(1784, 6).NumberFormat = "@"
$ws.Cells.Item(1784, 7).NumberFormat = "@"
$ws.Cells.Item(1784, 8).NumberFormat = "@"
$ws.Cells.Item(1784, 9).Value = 'select subString(''www.baidu.com'' from 5)'
$ws.Cells.Item(1784, 10).Value = 'baidu.com'
$ws.Cells.Item(1784, 11).Value = 'string_equals'

# Row 1785
$ws.Cells.Item(1785, 1).Value = 'dqlc1_1784'
$ws.Cells.Item(1785, 2).Value = 'y'
$ws.Cells.Item(1785, 3).Value = 'substring函数从指定起始索引处截取到最后'
$ws.Cells.Item(1785, 4).Value = 'String function'
$ws.Cells.Item(1785, 5).Value = 'subString'
$ws.Cells.Item(1785, 6).NumberFormat = "@"
$ws.Cells.Item(1785, 7).NumberFormat = "@"
$ws.Cells.Item(1785, 8).NumberFormat = "@"
$ws.Cells.Item(1785, 9).Value = 'select subString(''www.baidu.com'',5)'
$ws.Cells.Item(1785, 10).Value = 'baidu.com'
$ws.Cells.Item(1785, 11).Value = 'string_equals'

# Row 1786
$ws.Cells.Item(1786, 1).Value = 'dqlc1_1785'
$ws.Cells.Item(1786, 2).Value = 'y'
$ws.Cells.Item(1786, 3).Value = 'substring函数从指定起始索引处截取到最后'
$ws.Cells.Item(1786, 4).Value = 'String function'
$ws.Cells.Item(1786, 5).Value = 'subString'
$ws.Cells.Item(1786, 6).NumberFormat = "@"
$ws.Cells.Item(1786, 7).NumberFormat = "@"
$ws.Cells.Item(1786, 8).NumberFormat = "@"
$ws.Cells.Item(1786, 9).Value = 'select subString(''www.baidu.com'',1)'
$ws.Cells.Item(1786, 10).Value = 'www.baidu.com'
$ws.Cells.Item(1786, 11).Value = 'string_equals'

# Row 1787
$ws.Cells.Item(1787, 1).Value = 'dqlc1_1786'
$ws.Cells.Item(1787, 2).Value = 'y'
$ws.Cells.Item(1787, 3).Value = 'pow函数第二个参数为小数'
$ws.Cells.Item(1787, 4).Value = 'Numeric function'
$ws.Cells.Item(1787, 5).Value = 'pow'
$ws.Cells.Item(1787, 6).NumberFormat = "@"
$ws.Cells.Item(1787, 7).NumberFormat = "@"
$ws.Cells.Item(1787, 8).NumberFormat = "@"
$ws.Cells.Item(1787, 9).Value = 'select pow(16,0.5)'
$ws.Cells.Item(1787, 10).Value = '4'
$ws.Cells.Item(1787, 11).Value = 'string_equals'

# Update the sheet view to reflect the new scroll position / selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1754
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D1761").Select()
